# STEP 8 + STEP 9. Final.
# Updates the DM air matrix values in rows 7-15 (columns G-O) per the
# recalculated STEP 8 / STEP 9 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry describes the target cell, whether it should hold a numeric
# value ("n") or a text value ("s"), and the literal value to write.
$updates = @(
    @{Cell="H7"; Type="n"; Value=129600},
    @{Cell="I7"; Type="n"; Value=103680},
    @{Cell="J7"; Type="n"; Value=120960},
    @{Cell="K7"; Type="n"; Value=112320},
    @{Cell="L7"; Type="n"; Value=120960},
    @{Cell="M7"; Type="n"; Value=146880},
    @{Cell="N7"; Type="n"; Value=43200},
    @{Cell="O7"; Type="n"; Value=233280},

    @{Cell="G8"; Type="n"; Value=129600},
    @{Cell="I8"; Type="s"; Value="0"},
    @{Cell="J8"; Type="s"; Value="25920"},
    @{Cell="K8"; Type="s"; Value="17280"},
    @{Cell="L8"; Type="s"; Value="25920"},
    @{Cell="M8"; Type="s"; Value="0"},
    @{Cell="N8"; Type="s"; Value="0"},
    @{Cell="O8"; Type="s"; Value="60480"},

    @{Cell="G9"; Type="n"; Value=103680},
    @{Cell="H9"; Type="s"; Value="0"},
    @{Cell="J9"; Type="s"; Value="25920"},
    @{Cell="K9"; Type="s"; Value="17280"},
    @{Cell="L9"; Type="s"; Value="0"},
    @{Cell="M9"; Type="s"; Value="25920"},
    @{Cell="N9"; Type="s"; Value="17280"},
    @{Cell="O9"; Type="s"; Value="17280"},

    @{Cell="G10"; Type="n"; Value=120960},
    @{Cell="H10"; Type="s"; Value="25920"},
    @{Cell="I10"; Type="s"; Value="25920"},
    @{Cell="K10"; Type="s"; Value="8640"},
    @{Cell="L10"; Type="s"; Value="8640"},
    @{Cell="M10"; Type="s"; Value="17280"},
    @{Cell="N10"; Type="s"; Value="8640"},
    @{Cell="O10"; Type="s"; Value="25920"},

    @{Cell="G11"; Type="n"; Value=112320},
    @{Cell="H11"; Type="s"; Value="17280"},
    @{Cell="I11"; Type="s"; Value="17280"},
    @{Cell="J11"; Type="s"; Value="8640"},
    @{Cell="L11"; Type="s"; Value="17280"},
    @{Cell="M11"; Type="s"; Value="17280"},
    @{Cell="N11"; Type="s"; Value="8640"},
    @{Cell="O11"; Type="s"; Value="25920"},

    @{Cell="G12"; Type="n"; Value=120960},
    @{Cell="H12"; Type="s"; Value="25920"},
    @{Cell="I12"; Type="s"; Value="0"},
    @{Cell="J12"; Type="s"; Value="8640"},
    @{Cell="K12"; Type="s"; Value="17280"},
    @{Cell="M12"; Type="s"; Value="25920"},
    @{Cell="N12"; Type="s"; Value="8640"},
    @{Cell="O12"; Type="s"; Value="34560"},

    @{Cell="G13"; Type="n"; Value=146880},
    @{Cell="H13"; Type="s"; Value="0"},
    @{Cell="I13"; Type="s"; Value="25920"},
    @{Cell="J13"; Type="s"; Value="17280"},
    @{Cell="K13"; Type="s"; Value="17280"},
    @{Cell="L13"; Type="s"; Value="25920"},
    @{Cell="N13"; Type="s"; Value="0"},
    @{Cell="O13"; Type="s"; Value="60480"},

    @{Cell="G14"; Type="n"; Value=51840},
    @{Cell="H14"; Type="s"; Value=":"},
    @{Cell="I14"; Type="n"; Value=17280},
    @{Cell="J14"; Type="n"; Value=8640},
    @{Cell="K14"; Type="n"; Value=8640},
    @{Cell="L14"; Type="n"; Value=8640},
    @{Cell="M14"; Type="s"; Value=":"},
    @{Cell="O14"; Type="n"; Value=8640},

    @{Cell="G15"; Type="n"; Value=224640},
    @{Cell="H15"; Type="s"; Value="60480"},
    @{Cell="I15"; Type="s"; Value="17280"},
    @{Cell="J15"; Type="s"; Value="25920"},
    @{Cell="K15"; Type="s"; Value="25920"},
    @{Cell="L15"; Type="s"; Value="34560"},
    @{Cell="M15"; Type="s"; Value="60480"},
    @{Cell="N15"; Type="s"; Value="0"}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Type -eq "n") {
        $rng.Value = $u.Value
    } else {
        # Force text storage (even for numeric-looking strings like "0",
        # "25920", etc.) so the cell keeps its textual type, matching the
        # source data which stores these as strings rather than numbers.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    }
}
